# Weekly data update: insert a new price-report row for Orégano
# (Mercado Mayorista Lo Valledor de Santiago) as row 95, pushing the
# existing rows 95-111 down to 96-112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 95; everything below shifts down by one.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(95, 1).Value  = 6
$ws.Cells.Item(95, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(95, 3).Value  = "Metropolitana"
$ws.Cells.Item(95, 4).Value  = 44474
$ws.Cells.Item(95, 5).Value  = 13
$ws.Cells.Item(95, 6).Value  = 100112029
$ws.Cells.Item(95, 7).Value  = "Orégano"
$ws.Cells.Item(95, 8).Value  = "Sin especificar"
$ws.Cells.Item(95, 9).Value  = "Primera"
$ws.Cells.Item(95, 10).Value = 34
$ws.Cells.Item(95, 11).Value = 8500
$ws.Cells.Item(95, 12).Value = 9000
$ws.Cells.Item(95, 13).Value = 8735
$ws.Cells.Item(95, 14).Value = "$/docena de atados"
$ws.Cells.Item(95, 15).Value = "Región Metropolitana"
$ws.Cells.Item(95, 16).Value = 2912
$ws.Cells.Item(95, 17).Value = 3
$ws.Cells.Item(95, 18).Value = "Hortaliza"
